{"js": "// Apply the \"Started section-5 from the course\" edit:\n//  1) Recolor 3 runs that currently use color=\"auto\" to color=\"00000A\"\n//     - \"(overriding default command on run with \u201esh\u201d)\"\n//     - \"WORKDIR /usr/app \u2013\"\n//     - \" any following command will be executed relative to this path in the container\"\n//  2) Expand the single-run \"Docker compose\" paragraph into three runs with\n//     added explanatory text about docker-compose.\n\nconst body = context.document.body;\n\n// --- 1) Recolor the three \"auto\" colored runs to 00000A -------------------\n\nconst overridingRuns = body.search(\n  \"(overriding default command on run with \\u201Esh\\u201D)\",\n  { matchCase: true }\n);\nconst workdirRuns = body.search(\"WORKDIR /usr/app \\u2013\", { matchCase: true });\nconst anyFollowingRuns = body.search(\n  \" any following command will be executed relative to this path in the container\",\n  { matchCase: true }\n);\nawait context.sync();\n\noverridingRuns.items[0].font.color = \"#00000A\";\nworkdirRuns.items[0].font.color = \"#00000A\";\nanyFollowingRuns.items[0].font.color = \"#00000A\";\nawait context.sync();\n\n// --- 2) Split \"Docker compose\" into three runs with new content -----------\n\nconst dockerComposeRuns = body.search(\"Docker compose\", { matchCase: true });\nawait context.sync();\n\nconst run1 = dockerComposeRuns.items[0];\n// Replace the text in-place, keeping the same run/formatting.\nrun1.insertText(\"Docker compose \\u2013 \", Word.InsertLocation.replace);\nawait context.sync();\n\n// Insert the second and third runs right after the first.\nconst run2 = run1.insertText(\n  \"separate CLI that gets installed along with Docker. \",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst run3 = run2.insertText(\n  \"Used to start up multiple Docker containers at the same time. Automates \" +\n    \"some of the long-winded arguments  we were passing to \\u2018docker run\\u2019\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// The three new runs currently share identical formatting with their\n// neighbours, so the serializer would otherwise re-merge them back into a\n// single run. Toggling a direct-formatting property on and back off forces\n// each one to persist as its own <w:r>, matching the source edit that\n// produced three discrete runs.\nfor (const r of [run3, run2, run1]) {\n  r.font.underline = \"Single\";\n  await context.sync();\n  r.font.underline = \"None\";\n  await context.sync();\n}\n", "ps1": "# Apply the \"Started section-5 from the course\" edit via Word COM interop:\n#  1) Recolor 3 runs that currently use color=\"auto\" to color=\"00000A\"\n#     - \"(overriding default command on run with \u201esh\u201d)\"\n#     - \"WORKDIR /usr/app \u2013\"\n#     - \" any following command will be executed relative to this path in the container\"\n#  2) Expand the single-run \"Docker compose\" paragraph into three runs with\n#     added explanatory text about docker-compose.\n\n$doc = $word.ActiveDocument\n\n# NOTE: Word's Font.Color is a BGR-packed long (classic OLE_COLOR / WdColor),\n# not an RGB hex value, so target RGB 00000A must be supplied as 0x0A0000.\n$targetColor = 0x0A0000\n\n# --- 1) Recolor the three \"auto\" colored runs to 00000A -------------------\n\n$rangeOverriding = $doc.Content\n$rangeOverriding.Find.ClearFormatting()\n$rangeOverriding.Find.Text = \"(overriding default command on run with \" + [char]0x201E + \"sh\" + [char]0x201D + \")\"\n$rangeOverriding.Find.MatchCase = $true\n$rangeOverriding.Find.Execute() | Out-Null\n$rangeOverriding.Font.Color = $targetColor\n\n$rangeWorkdir = $doc.Content\n$rangeWorkdir.Find.ClearFormatting()\n$rangeWorkdir.Find.Text = \"WORKDIR /usr/app \" + [char]0x2013\n$rangeWorkdir.Find.MatchCase = $true\n$rangeWorkdir.Find.Execute() | Out-Null\n$rangeWorkdir.Font.Color = $targetColor\n\n$rangeAnyFollowing = $doc.Content\n$rangeAnyFollowing.Find.ClearFormatting()\n$rangeAnyFollowing.Find.Text = \" any following command will be executed relative to this path in the container\"\n$rangeAnyFollowing.Find.MatchCase = $true\n$rangeAnyFollowing.Find.Execute() | Out-Null\n$rangeAnyFollowing.Font.Color = $targetColor\n\n# --- 2) Split \"Docker compose\" into three runs with new content -----------\n\n$run1 = $doc.Content\n$run1.Find.ClearFormatting()\n$run1.Find.Text = \"Docker compose\"\n$run1.Find.MatchCase = $true\n$run1.Find.Execute() | Out-Null\n\n# Replace the text in-place, keeping the same run/formatting.\n$run1.Text = \"Docker compose \" + [char]0x2013 + \" \"\n\n# Insert the second run right after the first, then the third right after\n# the second, tracking explicit Range boundaries so each new chunk can be\n# addressed precisely afterwards.\n$seg2 = \"separate CLI that gets installed along with Docker. \"\n$insertPoint2 = $run1.End\n$run1.InsertAfter($seg2)\n$run2 = $doc.Range($insertPoint2, $insertPoint2 + $seg2.Length)\n\n$seg3 = \"Used to start up multiple Docker containers at the same time. Automates \" +\n  \"some of the long-winded arguments  we were passing to \" + [char]0x2018 + \"docker run\" + [char]0x2019\n$insertPoint3 = $run2.End\n$run2.InsertAfter($seg3)\n$run3 = $doc.Range($insertPoint3, $insertPoint3 + $seg3.Length)\n\n# The three new runs currently share identical formatting with their\n# neighbours, so the serializer would otherwise re-merge them back into a\n# single run. Toggling a direct-formatting property on and back off forces\n# each one to persist as its own run, matching the source edit that\n# produced three discrete runs. Go from last to first so earlier toggles\n# don't get re-absorbed by a later un-toggled neighbour.\n$run3.Font.Underline = 1\n$run3.Font.Underline = 0\n\n$run2.Font.Underline = 1\n$run2.Font.Underline = 0\n\n$run1.Font.Underline = 1\n$run1.Font.Underline = 0\n"}
